$wb = $excel.ActiveWorkbook

# --- Products sheet: Field Size column changes (VARCHAR max length 255 -> 128) ---
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Range("D4").Value2 = 128
$wsProducts.Range("D5").Value2 = 128
# D6 (producttype / FLOAT field) no longer has a Field Size value - clear it entirely
$wsProducts.Range("D6").Clear()
$wsProducts.Range("D7").Value2 = 128

# --- Users sheet: Field Size column changes (VARCHAR max length 255 -> 128) ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("D4").Value2 = 128
$wsUsers.Range("D5").Value2 = 128
$wsUsers.Range("D6").Value2 = 128
$wsUsers.Range("D7").Value2 = 128
$wsUsers.Range("D8").Value2 = 128
$wsUsers.Range("D9").Value2 = 128

# --- Update selections / active sheet ---
# Products tab loses focus, selection moves to D8
$wsProducts.Activate()
$wsProducts.Range("D8").Select()

# Users tab becomes the active/selected tab, selection moves to D9
$wsUsers.Activate()
$wsUsers.Range("D9").Select()
